# Update workbook for data as of 2022-05-08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (tab name) from "Through 2022-04-29" to "Through 2022-04-30"
$ws.Name = "Through 2022-04-30"

# Update the header label in I1 that shows the "through" date
$ws.Range("I1").Value = "2022 (through 04-30)"

# Update April row (row 5) current-year value
$ws.Range("I5").Value = 117

# Update Total row (row 14) current-year value
$ws.Range("I14").Value = 552
